# Update the "想去人数" (want-to-go count) figures in the F column across
# the affected sheets, matching the source data refresh captured in the
# commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 294
$ws1.Range("F18").Value = 412
$ws1.Range("F22").Value = 6868
$ws1.Range("F23").Value = 7415
$ws1.Range("F34").Value = 212
$ws1.Range("F35").Value = 268
$ws1.Range("F36").Value = 659
$ws1.Range("F38").Value = 1353

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F17").Value = 267

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2506

# Sheet "全部类型" (All Types) - aggregated view of the above sheets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 294
$ws4.Range("F23").Value = 412
$ws4.Range("F27").Value = 6868
$ws4.Range("F28").Value = 7415
$ws4.Range("F32").Value = 212
$ws4.Range("F34").Value = 268
$ws4.Range("F37").Value = 659
$ws4.Range("F42").Value = 1353
$ws4.Range("F49").Value = 267
